{"js": "// Rename the \"ReportPeriod:\" label to \"Rapport periode:\" in the exported\n// report header, matching the new table headers (kept in sync with the\n// report generator). Only the label text changes; the date range that\n// follows it is left untouched.\nconst oldLabel = \"ReportPeriod:\";\nconst newLabel = \"Rapport periode:\";\n\nconst body = context.document.body;\nconst results = body.search(oldLabel, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newLabel, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Rename the \"ReportPeriod:\" label to \"Rapport periode:\" in the exported\n# report header, matching the new table headers (kept in sync with the\n# report generator). Only the label text changes; the date range that\n# follows it is left untouched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"ReportPeriod:\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Rapport periode:\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, [ref]$null, [ref]$null, [ref]$find.Forward, [ref]$find.Wrap, [ref]$find.Format, [ref]$find.Replacement.Text, 2)\n"}
